# Weekly update: insert two new daily price rows (Alcachofa, Terminal
# Hortofrutícola Agro Chillán) ahead of the existing data, shifting the
# previously-recorded rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 21 (pushes old rows 21-46 down to 23-48).
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# New row 21
$ws.Cells.Item(21,1).Value  = 7
$ws.Cells.Item(21,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21,3).Value  = "Ñuble"
$ws.Cells.Item(21,4).Value  = 44789
$ws.Cells.Item(21,5).Value  = 16
$ws.Cells.Item(21,6).Value  = 100112013
$ws.Cells.Item(21,7).Value  = "Alcachofa"
$ws.Cells.Item(21,8).Value  = "Argentina(o)"
$ws.Cells.Item(21,9).Value  = "Primera"
$ws.Cells.Item(21,10).Value = 100
$ws.Cells.Item(21,11).Value = 14000
$ws.Cells.Item(21,12).Value = 15000
$ws.Cells.Item(21,13).Value = 14500
$ws.Cells.Item(21,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(21,15).Value = "Provincia de Limarí"
$ws.Cells.Item(21,16).Value = 290
$ws.Cells.Item(21,17).Value = 50
$ws.Cells.Item(21,18).Value = "Hortaliza"

# New row 22
$ws.Cells.Item(22,1).Value  = 7
$ws.Cells.Item(22,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(22,3).Value  = "Ñuble"
$ws.Cells.Item(22,4).Value  = 44789
$ws.Cells.Item(22,5).Value  = 16
$ws.Cells.Item(22,6).Value  = 100112013
$ws.Cells.Item(22,7).Value  = "Alcachofa"
$ws.Cells.Item(22,8).Value  = "Madrigal"
$ws.Cells.Item(22,9).Value  = "Primera"
$ws.Cells.Item(22,10).Value = 60
$ws.Cells.Item(22,11).Value = 13000
$ws.Cells.Item(22,12).Value = 14000
$ws.Cells.Item(22,13).Value = 13500
$ws.Cells.Item(22,14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(22,15).Value = "Provincia de Limarí"
$ws.Cells.Item(22,16).Value = 338
$ws.Cells.Item(22,17).Value = 40
$ws.Cells.Item(22,18).Value = "Hortaliza"
